$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the sheet-count ("খাতা/পত্রের সংখ্যা") values that were missing,
# which drives the per-row fee formulas in column I and the total in I32.
$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
